$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.878.32'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '1.706.29'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.68%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3950'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4092'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.524'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.007'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08830'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("E13").Value = '  +6.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001376'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.075'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.89%  '

$ws.Range("D17").Value = '1.703.53'
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07132'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.429'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.83%  '

$ws.Range("E22").Value = '  +1.11%  '

$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").Value = '24.876.91'
$ws.Range("E24").Value = '  +0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.071'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.350'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("E27").Value = '  +0.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.755'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.232'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.756'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.31%  '

$ws.Range("D33").Value = '1.891.93'
$ws.Range("E33").Value = '  +0.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08954'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.060'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02904'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.33%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09172'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7939'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.478'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7297'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.628'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.258'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.341'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '92.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.51%  '
